$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 43564.824
$ws.Range("I33").Value = 52720.316
$ws.Range("J33").Value = 76.25
$ws.Range("K33").Value = 52720.316
$ws.Range("L33").Value = 76.25
$ws.Range("M33").Value = -52491.316
$ws.Range("N33").Value = -534.25

# Row 64
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

# Row 112
$ws.Range("H112").Value = 2375
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2375
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 7125
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -9341

# Row 132
$ws.Range("H132").Value = 1799.8788
$ws.Range("I132").Value = 1942.8572
$ws.Range("J132").Value = 999.2
$ws.Range("K132").Value = 5828.571599999999
$ws.Range("L132").Value = 2997.6
$ws.Range("M132").Value = -3298.571599999999
$ws.Range("N132").Value = -8057.6

# Row 137
$ws.Range("H137").Value = 2300.7942
$ws.Range("I137").Value = 2219.8
$ws.Range("J137").Value = 2364.7368
$ws.Range("K137").Value = 6659.400000000001
$ws.Range("L137").Value = 7094.2104
$ws.Range("M137").Value = -4109.400000000001
$ws.Range("N137").Value = -12194.2104

# Row 138
$ws.Range("H138").Value = 2474.34
$ws.Range("I138").Value = 1662.9062
$ws.Range("J138").Value = 2856.1912
$ws.Range("K138").Value = 4988.7186
$ws.Range("L138").Value = 8568.5736
$ws.Range("M138").Value = 151.2813999999998
$ws.Range("N138").Value = -18848.5736

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1234.0851
$ws.Range("I61").Value = 1108.1765
$ws.Range("J61").Value = 1563.3846
$ws.Range("K61").Value = 1108.1765
$ws.Range("L61").Value = 1563.3846
$ws.Range("M61").Value = -896.1765
$ws.Range("N61").Value = -1987.3846

# Row 136
$ws.Range("H136").Value = 1234.0851
$ws.Range("I136").Value = 1108.1765
$ws.Range("J136").Value = 1563.3846
$ws.Range("K136").Value = 3324.5295
$ws.Range("L136").Value = 4690.1538
$ws.Range("M136").Value = -774.5295000000001
$ws.Range("N136").Value = -9790.1538

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 94225.55
$ws.Range("I134").Value = 4050.8
$ws.Range("J134").Value = 287457.16
$ws.Range("K134").Value = 12152.4
$ws.Range("L134").Value = 862371.48
$ws.Range("M134").Value = -9617.400000000001
$ws.Range("N134").Value = -867441.48

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1480
$ws.Range("I16").Value = 1471
$ws.Range("J16").Value = 1499.8
$ws.Range("K16").Value = 1471
$ws.Range("L16").Value = 1499.8
$ws.Range("M16").Value = -1184
$ws.Range("N16").Value = -2073.8

# Row 31
$ws.Range("H31").Value = 2973.3635
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 2973.3635
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 2973.3635
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -3563.3635

# Row 34
$ws.Range("H34").Value = 2973.3635
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 2973.3635
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 2973.3635
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -3377.3635

# Row 62
$ws.Range("H62").Value = 10453.846
$ws.Range("I62").Value = 11914.286
$ws.Range("J62").Value = 8750
$ws.Range("K62").Value = 11914.286
$ws.Range("L62").Value = 8750
$ws.Range("M62").Value = -11290.286
$ws.Range("N62").Value = -9998

# Row 65
$ws.Range("H65").Value = 10453.846
$ws.Range("I65").Value = 11914.286
$ws.Range("J65").Value = 8750
$ws.Range("K65").Value = 59571.43
$ws.Range("L65").Value = 43750
$ws.Range("M65").Value = -56451.43
$ws.Range("N65").Value = -49990

# Row 113
$ws.Range("H113").Value = 1480
$ws.Range("I113").Value = 1471
$ws.Range("J113").Value = 1499.8
$ws.Range("K113").Value = 1471
$ws.Range("L113").Value = 1499.8
$ws.Range("M113").Value = 699
$ws.Range("N113").Value = -5839.8

# Row 122
$ws.Range("H122").Value = 1250908.9
$ws.Range("I122").Value = 2000862.2
$ws.Range("K122").Value = 6002586.6
$ws.Range("M122").Value = -6000136.6

$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 5136.9565
$ws.Range("J39").Value = 5136.9565
$ws.Range("L39").Value = 15410.8695
$ws.Range("N39").Value = -15998.8695

# Row 68
$ws.Range("H68").Value = 1232.26
$ws.Range("I68").Value = 674.0755
$ws.Range("J68").Value = 1861.7021
$ws.Range("K68").Value = 2022.2265
$ws.Range("L68").Value = 5585.106299999999
$ws.Range("M68").Value = -1211.2265
$ws.Range("N68").Value = -7207.106299999999

# Row 71
$ws.Range("H71").Value = 1232.26
$ws.Range("I71").Value = 674.0755
$ws.Range("J71").Value = 1861.7021
$ws.Range("K71").Value = 6066.6795
$ws.Range("L71").Value = 16755.3189
$ws.Range("M71").Value = -2010.6795
$ws.Range("N71").Value = -24867.3189

# Row 126
$ws.Range("H126").Value = 4133.125
$ws.Range("I126").Value = 1638.3334
$ws.Range("J126").Value = 5630
$ws.Range("K126").Value = 4915.0002
$ws.Range("L126").Value = 16890
$ws.Range("M126").Value = 24.9997999999996
$ws.Range("N126").Value = -26770

$ws = $wb.Worksheets.Item("GSM")
# Row 136
$ws.Range("H136").Value = 18795.9
$ws.Range("J136").Value = 18795.9
$ws.Range("L136").Value = 56387.7
$ws.Range("N136").Value = -61487.7

# Row 138
$ws.Range("H138").Value = 40653.332
$ws.Range("J138").Value = 40653.332
$ws.Range("L138").Value = 40653.332
$ws.Range("N138").Value = -50933.332

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1188.8889
$ws.Range("I7").Value = 800
$ws.Range("J7").Value = 1300
$ws.Range("K7").Value = 800
$ws.Range("L7").Value = 1300
$ws.Range("M7").Value = -688
$ws.Range("N7").Value = -1524

# Row 50
$ws.Range("H50").Value = 6000
$ws.Range("J50").Value = 6000
$ws.Range("L50").Value = 6000
$ws.Range("N50").Value = -7274

# Row 54
$ws.Range("H54").Value = 9999.5
$ws.Range("J54").Value = 9999.5
$ws.Range("L54").Value = 9999.5
$ws.Range("N54").Value = -11287.5

# Row 61
$ws.Range("H61").Value = 1184.0526
$ws.Range("I61").Value = 1182.5454
$ws.Range("K61").Value = 1182.5454
$ws.Range("M61").Value = -980.5454

# Row 113
$ws.Range("H113").Value = 1184.0526
$ws.Range("I113").Value = 1182.5454
$ws.Range("K113").Value = 1182.5454
$ws.Range("M113").Value = 987.4546

# Row 126
$ws.Range("H126").Value = 1188.8889
$ws.Range("I126").Value = 800
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 2400
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = 70
$ws.Range("N126").Value = -8840

# Row 136
$ws.Range("H136").Value = 2569.7556
$ws.Range("I136").Value = 1634.6296
$ws.Range("J136").Value = 3972.4443
$ws.Range("K136").Value = 4903.8888
$ws.Range("L136").Value = 11917.3329
$ws.Range("M136").Value = -2353.8888
$ws.Range("N136").Value = -17017.3329

# Row 139
$ws.Range("H139").Value = 52500
$ws.Range("J139").Value = 52500
$ws.Range("L139").Value = 52500
$ws.Range("N139").Value = -62780

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2065.0334
$ws.Range("I132").Value = 794.1667
$ws.Range("J132").Value = 2912.2778
$ws.Range("K132").Value = 2382.5001
$ws.Range("L132").Value = 8736.8334
$ws.Range("M132").Value = 147.4998999999998
$ws.Range("N132").Value = -13796.8334

# Row 136
$ws.Range("H136").Value = 2574.6458
$ws.Range("I136").Value = 3194.913
$ws.Range("K136").Value = 9584.739
$ws.Range("M136").Value = -7034.739

# Row 138
$ws.Range("H138").Value = 48280
$ws.Range("J138").Value = 48280
$ws.Range("L138").Value = 48280
$ws.Range("N138").Value = -58560
